# Paso de usarse un array de niveles a usarse una tabla hash de niveles.
# Applies the cell changes from the "Nivel Facil" worksheet (the active sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# Row 1: move "p" from D1 to C1, and "d" from H1 to I1
$ws.Range("D1").Clear()
$ws.Range("C1").Value = "p"

$ws.Range("H1").Clear()
$ws.Range("I1").Value = "d"

# Row 4: add new "cv" entry
$ws.Range("B4").Value = "cv"

# Row 5: G5 ("f") is removed
$ws.Range("G5").Clear()

# Row 7: add new "cp" entry
$ws.Range("E7").Value = "cp"

# Row 9: E9 changes from "c" to "cv"
$ws.Range("E9").Value = "cv"

# Update the active selection to E3
[void]$ws.Range("E3").Select()
